# Weekly update for "Comercializadora del Agro de Limarí - Alcachofa".
# The previous week's rows (142-145) are overwritten with the new week's
# figures, a new 4th variety row is appended (new row 145), and the three
# rows of data that were displaced (old rows 143-145) are appended as new
# rows 146-148, preserving history.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 142: Argentina(o) / Primera -> new week figures ---
$ws.Range("D142").Value = 44461
$ws.Range("J142").Value = 1000
$ws.Range("K142").Value = 7000
$ws.Range("L142").Value = 8000
$ws.Range("M142").Value = 7500
$ws.Range("P142").Value = 150

# --- Row 143: Española / Primera -> new week figures ---
$ws.Range("D143").Value = 44461
$ws.Range("J143").Value = 1200
$ws.Range("K143").Value = 8000
$ws.Range("L143").Value = 10000
$ws.Range("M143").Value = 9000
$ws.Range("P143").Value = 300

# --- Row 144: becomes Madrigal / Primera for the new week ---
$ws.Range("D144").Value = 44461
$ws.Range("H144").Value = "Madrigal"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 1200
$ws.Range("K144").Value = 6000
$ws.Range("L144").Value = 7000
$ws.Range("M144").Value = 6500
$ws.Range("P144").Value = 162

# --- Row 145: becomes Argentina(o) / Primera for the new week ---
$ws.Range("H145").Value = "Argentina(o)"
$ws.Range("J145").Value = 700
$ws.Range("N145").Value = "$/caja 50 unidades"
$ws.Range("P145").Value = 270
$ws.Range("Q145").Value = 50

# --- New row 146: Española / Primera (previous week data, carried down) ---
$ws.Range("A146").Value = 2
$ws.Range("B146").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = 44357
$ws.Range("D146").NumberFormat = $ws.Range("D140").NumberFormat
$ws.Range("E146").Value = 4
$ws.Range("F146").Value = 100112013
$ws.Range("G146").Value = "Alcachofa"
$ws.Range("H146").Value = "Española"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 800
$ws.Range("K146").Value = 14000
$ws.Range("L146").Value = 15000
$ws.Range("M146").Value = 14500
$ws.Range("N146").Value = "$/caja 30 unidades"
$ws.Range("O146").Value = "Provincia de Limarí"
$ws.Range("P146").Value = 483
$ws.Range("Q146").Value = 30
$ws.Range("R146").Value = "Hortaliza"

# --- New row 147: Española / Segunda (previous week data, carried down) ---
$ws.Range("A147").Value = 2
$ws.Range("B147").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = 44357
$ws.Range("D147").NumberFormat = $ws.Range("D140").NumberFormat
$ws.Range("E147").Value = 4
$ws.Range("F147").Value = 100112013
$ws.Range("G147").Value = "Alcachofa"
$ws.Range("H147").Value = "Española"
$ws.Range("I147").Value = "Segunda"
$ws.Range("J147").Value = 500
$ws.Range("K147").Value = 12000
$ws.Range("L147").Value = 13000
$ws.Range("M147").Value = 12500
$ws.Range("N147").Value = "$/caja 40 unidades"
$ws.Range("O147").Value = "Provincia de Limarí"
$ws.Range("P147").Value = 312
$ws.Range("Q147").Value = 40
$ws.Range("R147").Value = "Hortaliza"

# --- New row 148: Madrigal / Primera (previous week data, carried down) ---
$ws.Range("A148").Value = 2
$ws.Range("B148").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44357
$ws.Range("D148").NumberFormat = $ws.Range("D140").NumberFormat
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = 100112013
$ws.Range("G148").Value = "Alcachofa"
$ws.Range("H148").Value = "Madrigal"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 500
$ws.Range("K148").Value = 13000
$ws.Range("L148").Value = 14000
$ws.Range("M148").Value = 13500
$ws.Range("N148").Value = "$/caja 40 unidades"
$ws.Range("O148").Value = "Provincia de Limarí"
$ws.Range("P148").Value = 338
$ws.Range("Q148").Value = 40
$ws.Range("R148").Value = "Hortaliza"
